# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 58829616
$ws.Range("I132").Value = 71435320
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 214305960
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -214303430
$ws.Range("N132").Value = -14057
$ws.Range("H135").Value = 2130.3333
$ws.Range("I135").Value = 2130.3333
$ws.Range("K135").Value = 19172.9997
$ws.Range("M135").Value = -16637.9997
$ws.Range("H137").Value = 3840.7058
$ws.Range("I137").Value = 4837.125
$ws.Range("K137").Value = 14511.375
$ws.Range("M137").Value = -11961.375
$ws.Range("H138").Value = 1533.5555
$ws.Range("I138").Value = 872.6667
$ws.Range("J138").Value = 1680.4198
$ws.Range("K138").Value = 2618.0001
$ws.Range("L138").Value = 5041.2594
$ws.Range("M138").Value = 2521.9999
$ws.Range("N138").Value = -15321.2594

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10528.167
$ws.Range("I32").Value = 10663
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 10663
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -10376
$ws.Range("N32").Value = -5574
$ws.Range("H45").Value = 3623.7273
$ws.Range("J45").Value = 5202.75
$ws.Range("L45").Value = 5202.75
$ws.Range("N45").Value = -5956.75
$ws.Range("H61").Value = 987432
$ws.Range("I61").Value = 1617130.1
$ws.Range("K61").Value = 1617130.1
$ws.Range("M61").Value = -1616918.1
$ws.Range("H74").Value = 6216711.5
$ws.Range("I74").Value = 23810922
$ws.Range("J74").Value = 6990.0586
$ws.Range("K74").Value = 23810922
$ws.Range("L74").Value = 6990.0586
$ws.Range("M74").Value = -23810048
$ws.Range("N74").Value = -8738.0586
$ws.Range("H77").Value = 6216711.5
$ws.Range("I77").Value = 23810922
$ws.Range("J77").Value = 6990.0586
$ws.Range("K77").Value = 119054610
$ws.Range("L77").Value = 34950.29300000001
$ws.Range("M77").Value = -119050242
$ws.Range("N77").Value = -43686.29300000001
$ws.Range("H132").Value = 2851174.8
$ws.Range("I132").Value = 3346557.5
$ws.Range("J132").Value = 2725
$ws.Range("K132").Value = 10039672.5
$ws.Range("L132").Value = 8175
$ws.Range("M132").Value = -10037142.5
$ws.Range("N132").Value = -13235
$ws.Range("H136").Value = 987432
$ws.Range("I136").Value = 1617130.1
$ws.Range("K136").Value = 4851390.300000001
$ws.Range("M136").Value = -4848840.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3514672.2
$ws.Range("I86").Value = 6064224.5
$ws.Range("J86").Value = 9037.875
$ws.Range("K86").Value = 6064224.5
$ws.Range("L86").Value = 9037.875
$ws.Range("M86").Value = -6063101.5
$ws.Range("N86").Value = -11283.875
$ws.Range("H89").Value = 3514672.2
$ws.Range("I89").Value = 6064224.5
$ws.Range("J89").Value = 9037.875
$ws.Range("K89").Value = 30321122.5
$ws.Range("L89").Value = 45189.375
$ws.Range("M89").Value = -30315506.5
$ws.Range("N89").Value = -56421.375
$ws.Range("H94").Value = 5927
$ws.Range("I94").Value = 2289.2
$ws.Range("J94").Value = 13202.6
$ws.Range("K94").Value = 2289.2
$ws.Range("L94").Value = 13202.6
$ws.Range("M94").Value = -1838.2
$ws.Range("N94").Value = -14104.6
$ws.Range("H118").Value = 39998.668
$ws.Range("J118").Value = 39998.668
$ws.Range("L118").Value = 39998.668
$ws.Range("N118").Value = -43312.668
$ws.Range("H134").Value = 5356.547
$ws.Range("I134").Value = 1767.9762
$ws.Range("J134").Value = 9923.817999999999
$ws.Range("K134").Value = 5303.9286
$ws.Range("L134").Value = 29771.454
$ws.Range("M134").Value = -2768.9286
$ws.Range("N134").Value = -34841.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2743.9473
$ws.Range("I31").Value = 2539.8125
$ws.Range("K31").Value = 2539.8125
$ws.Range("M31").Value = -2244.8125
$ws.Range("H34").Value = 2743.9473
$ws.Range("I34").Value = 2539.8125
$ws.Range("K34").Value = 2539.8125
$ws.Range("M34").Value = -2337.8125
$ws.Range("H86").Value = 4952.7
$ws.Range("I86").Value = 4940.875
$ws.Range("K86").Value = 4940.875
$ws.Range("M86").Value = -3817.875
$ws.Range("H89").Value = 4952.7
$ws.Range("I89").Value = 4940.875
$ws.Range("K89").Value = 24704.375
$ws.Range("M89").Value = -19088.375
$ws.Range("H99").Value = 3197.5
$ws.Range("I99").Value = 2799
$ws.Range("K99").Value = 2799
$ws.Range("M99").Value = -1301
$ws.Range("H126").Value = 3197.5
$ws.Range("I126").Value = 2799
$ws.Range("K126").Value = 8397
$ws.Range("M126").Value = -5927

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 750
$ws.Range("I47").Value = 1000
$ws.Range("J47").Value = 500
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 1500
$ws.Range("M47").Value = -2569
$ws.Range("N47").Value = -2362
$ws.Range("H59").Value = 5285.7144
$ws.Range("I59").Value = 3333.3333
$ws.Range("J59").Value = 6750
$ws.Range("K59").Value = 9999.999899999999
$ws.Range("L59").Value = 20250
$ws.Range("M59").Value = -9459.999899999999
$ws.Range("N59").Value = -21330
$ws.Range("H92").Value = 806.38464
$ws.Range("J92").Value = 842.5454999999999
$ws.Range("L92").Value = 2527.6365
$ws.Range("N92").Value = -5023.6365
$ws.Range("H98").Value = 2074.7
$ws.Range("I98").Value = 786.75
$ws.Range("J98").Value = 2933.3333
$ws.Range("K98").Value = 2360.25
$ws.Range("L98").Value = 8799.999899999999
$ws.Range("M98").Value = -862.25
$ws.Range("N98").Value = -11795.9999
$ws.Range("H132").Value = 2724.3171
$ws.Range("I132").Value = 1330.1428
$ws.Range("J132").Value = 3447.2222
$ws.Range("K132").Value = 11971.2852
$ws.Range("L132").Value = 31024.9998
$ws.Range("M132").Value = -9441.2852
$ws.Range("N132").Value = -36084.99980000001
$ws.Range("H134").Value = 4768.385
$ws.Range("I134").Value = 3401.625
$ws.Range("J134").Value = 6955.2
$ws.Range("K134").Value = 10204.875
$ws.Range("L134").Value = 20865.6
$ws.Range("M134").Value = -5134.875
$ws.Range("N134").Value = -31005.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 10005502
$ws.Range("I14").Value = 10005502
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 10005502
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -10005334
$ws.Range("N14").ClearContents()
$ws.Range("H39").Value = 97630
$ws.Range("J39").Value = 97630
$ws.Range("L39").Value = 97630
$ws.Range("N39").Value = -98694
$ws.Range("H80").Value = 3252.0588
$ws.Range("I80").Value = 3126.5715
$ws.Range("J80").Value = 3339.9
$ws.Range("K80").Value = 3126.5715
$ws.Range("L80").Value = 3339.9
$ws.Range("M80").Value = -2128.5715
$ws.Range("N80").Value = -5335.9
$ws.Range("H83").Value = 3252.0588
$ws.Range("I83").Value = 3126.5715
$ws.Range("J83").Value = 3339.9
$ws.Range("K83").Value = 15632.8575
$ws.Range("L83").Value = 16699.5
$ws.Range("M83").Value = -10640.8575
$ws.Range("N83").Value = -26683.5
$ws.Range("H97").Value = 1492.5
$ws.Range("I97").Value = 551.5238000000001
$ws.Range("J97").Value = 2654.8823
$ws.Range("K97").Value = 551.5238000000001
$ws.Range("L97").Value = 2654.8823
$ws.Range("M97").Value = -55.52380000000005
$ws.Range("N97").Value = -3646.8823
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 29493.9
$ws.Range("J136").Value = 29493.9
$ws.Range("L136").Value = 88481.70000000001
$ws.Range("N136").Value = -93581.70000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5485.1562
$ws.Range("I40").Value = 5445.08
$ws.Range("J40").Value = 5628.2856
$ws.Range("K40").Value = 5445.08
$ws.Range("L40").Value = 5628.2856
$ws.Range("M40").Value = -5309.08
$ws.Range("N40").Value = -5900.2856
$ws.Range("H68").Value = 2766
$ws.Range("I68").Value = 2778.6
$ws.Range("J68").Value = 2703
$ws.Range("K68").Value = 2778.6
$ws.Range("L68").Value = 2703
$ws.Range("M68").Value = -2029.6
$ws.Range("N68").Value = -4201
$ws.Range("H71").Value = 2766
$ws.Range("I71").Value = 2778.6
$ws.Range("J71").Value = 2703
$ws.Range("K71").Value = 13893
$ws.Range("L71").Value = 13515
$ws.Range("M71").Value = -10149
$ws.Range("N71").Value = -21003
$ws.Range("H140").Value = 91999.164
$ws.Range("J140").Value = 91999.164
$ws.Range("L140").Value = 91999.164
$ws.Range("N140").Value = -102359.164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 41040.875
$ws.Range("I2").Value = 44047
$ws.Range("J2").Value = 19998
$ws.Range("K2").Value = 44047
$ws.Range("L2").Value = 19998
$ws.Range("M2").Value = -43935
$ws.Range("N2").Value = -20222
$ws.Range("H81").Value = 55588.35
$ws.Range("I81").Value = 113653.336
$ws.Range("J81").Value = 8080.636
$ws.Range("K81").Value = 227306.672
$ws.Range("L81").Value = 16161.272
$ws.Range("M81").Value = -226245.672
$ws.Range("N81").Value = -18283.272
$ws.Range("H84").Value = 55588.35
$ws.Range("I84").Value = 113653.336
$ws.Range("J84").Value = 8080.636
$ws.Range("K84").Value = 1136533.36
$ws.Range("L84").Value = 80806.36
$ws.Range("M84").Value = -1131229.36
$ws.Range("N84").Value = -91414.36
$ws.Range("H113").Value = 4906635.5
$ws.Range("I113").Value = 5560647
$ws.Range("J113").Value = 1550
$ws.Range("K113").Value = 16681941
$ws.Range("L113").Value = 4650
$ws.Range("M113").Value = -16679771
$ws.Range("N113").Value = -8990
$ws.Range("H122").Value = 2871.1428
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 999.6667
$ws.Range("I132").Value = 999.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2999.0001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -469.0001000000002
$ws.Range("N132").ClearContents()
